# Add a new "localdb" command-type column/category to the "#system" lookup
# sheet, plus a new defined name, to support the new local database
# automation commands (cloneTable, dropTables, exportCSV, importRecords,
# purge, runSQLs).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# ---------------------------------------------------------------------------
# Step 1: insert a brand-new column at "N" (the existing "macro" column and
# everything to its right - mail, number, pdf, rdbms, redis, sms, sound,
# ssh, step, web, webalert, webcookie, ws, ws.async, xml - shift one column
# to the right, from N:AC to O:AD).
# ---------------------------------------------------------------------------
$ws.Columns("N").Insert()

# ---------------------------------------------------------------------------
# Step 2: populate the freshly inserted column N with the "localdb" header
# and its six function names.
# ---------------------------------------------------------------------------
$localdbValues = @(
    "localdb",
    "cloneTable(var,source,target)",
    "dropTables(var,tables)",
    "exportCSV(sql,output)",
    "importRecords(var,sourceDb,sql,table)",
    "purge(var)",
    "runSQLs(var,sqls)"
)
for ($i = 0; $i -lt $localdbValues.Length; $i++) {
    $ws.Cells.Item($i + 1, 14).Value2 = $localdbValues[$i]
}

# ---------------------------------------------------------------------------
# Step 3: insert "localdb" into the alphabetical command-type list held in
# column A, between "json" (row 13) and "macro" (row 14). Only column A
# shifts down by one row (A14:A29 -> A15:A30); every other column keeps its
# existing row alignment, so this is done with a manual value shift rather
# than a full row/range insert.
# ---------------------------------------------------------------------------
for ($r = 29; $r -ge 14; $r--) {
    $ws.Cells.Item($r + 1, 1).Value2 = $ws.Cells.Item($r, 1).Value2
}
$ws.Cells.Item(14, 1).Value2 = "localdb"

# ---------------------------------------------------------------------------
# Step 4: fix up all the defined names whose ranges reference the columns
# that just shifted right by one letter, and register the new "localdb"
# and "target" ranges.
# ---------------------------------------------------------------------------
$wb.Names.Item("macro").RefersTo      = "='#system'!`$O`$2:`$O`$4"
$wb.Names.Item("mail").RefersTo       = "='#system'!`$P`$2:`$P`$2"
$wb.Names.Item("number").RefersTo     = "='#system'!`$Q`$2:`$Q`$16"
$wb.Names.Item("pdf").RefersTo        = "='#system'!`$R`$2:`$R`$16"
$wb.Names.Item("rdbms").RefersTo      = "='#system'!`$S`$2:`$S`$7"
$wb.Names.Item("redis").RefersTo      = "='#system'!`$T`$2:`$T`$10"
$wb.Names.Item("sms").RefersTo        = "='#system'!`$U`$2:`$U`$2"
$wb.Names.Item("sound").RefersTo      = "='#system'!`$V`$2:`$V`$5"
$wb.Names.Item("ssh").RefersTo        = "='#system'!`$W`$2:`$W`$9"
$wb.Names.Item("step").RefersTo       = "='#system'!`$X`$2:`$X`$4"
$wb.Names.Item("target").RefersTo     = "='#system'!`$A`$2:`$A`$30"
$wb.Names.Item("web").RefersTo        = "='#system'!`$Y`$2:`$Y`$127"
$wb.Names.Item("webalert").RefersTo   = "='#system'!`$Z`$2:`$Z`$8"
$wb.Names.Item("webcookie").RefersTo  = "='#system'!`$AA`$2:`$AA`$8"
$wb.Names.Item("ws").RefersTo         = "='#system'!`$AB`$2:`$AB`$17"
$wb.Names.Item("ws.async").RefersTo   = "='#system'!`$AC`$2:`$AC`$8"
$wb.Names.Item("xml").RefersTo        = "='#system'!`$AD`$2:`$AD`$21"

$wb.Names.Add("localdb", "='#system'!`$N`$2:`$N`$7")
